$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '81.177.92'
$ws.Range("E2").Value = '  +5.01%  '

$ws.Range("D3").Value = '3.183.64'
$ws.Range("E3").Value = '  +1.38%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '209.76'
$ws.Range("E5").Value = '  +3.63%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '636.02'
$ws.Range("E6").Value = '  +1.06%  '

$ws.Range("E7").Value = '  +27.84%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.592'
$ws.Range("E9").Value = '  +4.10%  '

$ws.Range("D10").Value = '3.179.07'
$ws.Range("E10").Value = '  +1.29%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.594'
$ws.Range("E11").Value = '  +10.39%  '

$ws.Range("E12").Value = '  +19.26%  '

$ws.Range("E13").Value = '  +2.30%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.39'
$ws.Range("E14").Value = '  -1.07%  '

$ws.Range("D15").Value = '3.765.93'
$ws.Range("E15").Value = '  +1.32%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '32.19'
$ws.Range("E16").Value = '  +5.89%  '

$ws.Range("D17").Value = '81.116.11'
$ws.Range("E17").Value = '  +5.05%  '

$ws.Range("D18").Value = '3.172.89'
$ws.Range("E18").Value = '  +1.17%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.25'
$ws.Range("E19").Value = '  +15.08%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.37'
$ws.Range("E20").Value = '  +2.92%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.26'
$ws.Range("E21").Value = '  +0.78%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '441.36'
$ws.Range("E22").Value = '  +1.59%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.24'
$ws.Range("E23").Value = '  +10.53%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.08'
$ws.Range("E24").Value = '  +5.03%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.09'
$ws.Range("E25").Value = '  +10.09%  '

$ws.Range("E26").Value = '  +6.93%  '

$ws.Range("D27").Value = '3.348.55'
$ws.Range("E27").Value = '  +1.36%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '77.00'
$ws.Range("E28").Value = '  +2.03%  '

$ws.Range("E29").Value = '  +13.38%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.19%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.22'
$ws.Range("E31").Value = '  +6.03%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.997'
$ws.Range("E32").Value = '  -0.31%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '572.92'
$ws.Range("E33").Value = '  +10.58%  '

$ws.Range("E34").Value = '  +2.98%  '

$ws.Range("B35").Value = 'Cronos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.143'
$ws.Range("E35").Value = '  +33.69%  '

$ws.Range("E36").Value = '  +4.39%  '

$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.152'
$ws.Range("E37").Value = '  +12.94%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '23.29'
$ws.Range("E38").Value = '  +4.16%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.998'
$ws.Range("E39").Value = '  -0.18%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.415'
$ws.Range("E40").Value = '  +5.65%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.12'
$ws.Range("E41").Value = '  +24.85%  '

$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.07'
$ws.Range("E42").Value = '  +18.40%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.98'
$ws.Range("E43").Value = '  +11.71%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '20.79'
$ws.Range("E44").Value = '  +3.60%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '160.22'
$ws.Range("E45").Value = '  -2.18%  '

$ws.Range("E46").Value = '  -0.04%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '188.90'
$ws.Range("E47").Value = '  -3.09%  '

$ws.Range("B48").Value = 'OKB'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '45.20'
$ws.Range("E48").Value = '  +6.24%  '

$ws.Range("B49").Value = 'ImmutableX'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.35'
$ws.Range("E49").Value = '  +6.01%  '

$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.780'
$ws.Range("E50").Value = '  -1.88%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.31'
$ws.Range("E51").Value = '  +5.84%  '
